$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 54.86376272656823

$ws.Range("N2:N6").Value = $newValue
